# ajuste: corrigindo as categorias
# Adds a "Total" column (row-wise sums) and two new rows: "Outros" (a
# residual / other-causes category) and "Total" (column-wise sums across
# the disease categories).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Total" column header -------------------------------------------------
$ws.Range("T1").Value = "Total"

# --- Row-wise totals for the existing disease rows (2-6) -----------------------
$ws.Range("T2").Value = 84467
$ws.Range("T3").Value = 12093
$ws.Range("T4").Value = 40409
$ws.Range("T5").Value = 17974
$ws.Range("T6").Value = 55706

# --- New row 7: "Outros" (Other causes) -----------------------------------------
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 5694
$ws.Range("C7").Value = 212
$ws.Range("D7").Value = 341
$ws.Range("E7").Value = 1718
$ws.Range("F7").Value = 3045
$ws.Range("G7").Value = 2972
$ws.Range("H7").Value = 3672
$ws.Range("I7").Value = 4516
$ws.Range("J7").Value = 5366
$ws.Range("K7").Value = 6186
$ws.Range("L7").Value = 7660
$ws.Range("M7").Value = 9494
$ws.Range("N7").Value = 11559
$ws.Range("O7").Value = 12841
$ws.Range("P7").Value = 13078
$ws.Range("Q7").Value = 12987
$ws.Range("R7").Value = 37264
$ws.Range("S7").Value = 381
$ws.Range("T7").Value = 138986

# --- New row 8: "Total" (sum across all disease categories, rows 2-7) ----------
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 6141
$ws.Range("C8").Value = 364
$ws.Range("D8").Value = 544
$ws.Range("E8").Value = 2185
$ws.Range("F8").Value = 3858
$ws.Range("G8").Value = 4060
$ws.Range("H8").Value = 5423
$ws.Range("I8").Value = 7337
$ws.Range("J8").Value = 9811
$ws.Range("K8").Value = 12822
$ws.Range("L8").Value = 17674
$ws.Range("M8").Value = 24309
$ws.Range("N8").Value = 31874
$ws.Range("O8").Value = 36854
$ws.Range("P8").Value = 38682
$ws.Range("Q8").Value = 38869
$ws.Range("R8").Value = 108367
$ws.Range("S8").Value = 461
$ws.Range("T8").Value = 349635
